$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lookups")

# Insert a new column before column F (6) so that every existing column from
# F (Gender) onward shifts one column to the right. This makes room for a new
# "cuts" column that holds the Grade/Region menu-item values (one entry per
# config), which is what lets the UI build a separate menu for each config.
$ws.Columns.Item(6).Insert()

# Populate the newly inserted column F with the cuts values (mirrors the
# Grade/Region values already present in columns A/C).
$ws.Range("F1").Value = "Grade"
$ws.Range("F2").Value = "Region"
$ws.Range("F3").Value = "Region"

# Give the new column the same width as its neighbours.
$ws.Columns.Item(6).ColumnWidth = 9.1

# Update the defined names that pointed at the old layout so they reflect the
# inserted column. cuts_config keeps referring to the same fixed range.
$wb.Names.Item("cuts").RefersTo = "=Lookups!`$F`$1:`$F`$2"
$wb.Names.Add("cuts_historical", "=Lookups!`$F`$3:`$F`$3")
$wb.Names.Item("default_menu").RefersTo = "=Lookups!`$F`$2:`$F`$101"
$wb.Names.Item("default_mapping").RefersTo = "=Lookups!`$F`$2:`$G`$101"
$wb.Names.Item("default_menu_start").RefersTo = "=Lookups!`$F`$2"
$wb.Names.Item("cuts_head").RefersTo = "=Lookups!`$G`$1:`$U`$1"
$wb.Names.Item("zero_string").RefersTo = "=Lookups!`$V`$1"
